$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2 and 3 with new computed values
$ws.Range("B2").Value = 7.764885014274649
$ws.Range("C2").Value = 0.07000349873623612

$ws.Range("B3").Value = 6.1607253821163
$ws.Range("C3").Value = 0.07088297468101119

# Add new rows 4-6 following the same pattern as rows 2-3
$ws.Range("A4").Value = 900
$ws.Range("B4").Value = 3.318924097591853
$ws.Range("C4").Value = 0.07433659649167967
$ws.Range("D4").Value = "/"

$ws.Range("A5").Value = 1250
$ws.Range("B5").Value = 2.334301487175983
$ws.Range("C5").Value = 0.07810766883554814
$ws.Range("D5").Value = "/"

$ws.Range("A6").Value = 1500
$ws.Range("B6").Value = 1.916840969215296
$ws.Range("C6").Value = 0.08151118187122558
$ws.Range("D6").Value = "/"

# Apply the same formatting as other A-column cells (bold, thin border, centered) to new A cells
$ws.Range("A3").Copy()
$ws.Range("A4:A6").PasteSpecial(-4122)
